# Apply updated "dSF" (column F) values for a set of rows on Sheet1.
# These reflect a repull/recalculation of the data (per commit message:
# "repull data, push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 1
    4  = 1
    5  = 1
    7  = -4
    8  = 0
    10 = -1
    11 = -1
    12 = 1
    16 = 1
    17 = 1
    28 = -1
    32 = 0
    35 = 0
    44 = -1
    46 = 3
    50 = 2
    54 = 0
    58 = 5
    66 = -1
    68 = -5
    69 = 0
    70 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
